# Auto-generated script to apply scheduled market-data refresh to Marilith_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H-N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2641.0833
$ws.Range("J17").Value = 2641.0833
$ws.Range("L17").Value = 7923.249899999999
$ws.Range("N17").Value = -8259.249899999999
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5350
$ws.Range("H42").Value = 48.75
$ws.Range("I42").Value = 41.666668
$ws.Range("K42").Value = 125.000004
$ws.Range("M42").Value = 104.999996
$ws.Range("H53").Value = 246
$ws.Range("I53").Value = 331.625
$ws.Range("K53").Value = 331.625
$ws.Range("M53").Value = 305.375
$ws.Range("H92").Value = 1563.3334
$ws.Range("I92").Value = 1665
$ws.Range("J92").Value = 750
$ws.Range("K92").Value = 1665
$ws.Range("L92").Value = 750
$ws.Range("M92").Value = -417
$ws.Range("N92").Value = -3246
$ws.Range("H100").Value = 2899.4443
$ws.Range("I100").Value = 2800
$ws.Range("K100").Value = 2800
$ws.Range("M100").Value = -2259
$ws.Range("H111").Value = 195
$ws.Range("I111").Value = 195
$ws.Range("K111").Value = 585
$ws.Range("M111").Value = 2482
$ws.Range("H125").Value = 4090.6667
$ws.Range("I125").Value = 3188.1428
$ws.Range("J125").Value = 7249.5
$ws.Range("K125").Value = 28693.2852
$ws.Range("L125").Value = 65245.5
$ws.Range("M125").Value = -26233.2852
$ws.Range("N125").Value = -70165.5
$ws.Range("H132").Value = 2206.6538
$ws.Range("I132").Value = 2094.92
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6284.76
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3754.76
$ws.Range("N132").Value = -20060
$ws.Range("H138").Value = 2561.7
$ws.Range("I138").Value = 1420
$ws.Range("K138").Value = 4260
$ws.Range("M138").Value = 880

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2374.8
$ws.Range("I74").Value = 1821.1428
$ws.Range("J74").Value = 3666.6667
$ws.Range("K74").Value = 1821.1428
$ws.Range("L74").Value = 3666.6667
$ws.Range("M74").Value = -947.1428000000001
$ws.Range("N74").Value = -5414.6667
$ws.Range("H77").Value = 2374.8
$ws.Range("I77").Value = 1821.1428
$ws.Range("J77").Value = 3666.6667
$ws.Range("K77").Value = 9105.714
$ws.Range("L77").Value = 18333.3335
$ws.Range("M77").Value = -4737.714
$ws.Range("N77").Value = -27069.3335
$ws.Range("H97").Value = 1076
$ws.Range("I97").Value = 917
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 917
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -421
$ws.Range("N97").Value = -2492
$ws.Range("H132").Value = 828.12
$ws.Range("I132").Value = 845.1739
$ws.Range("J132").Value = 632
$ws.Range("K132").Value = 2535.5217
$ws.Range("L132").Value = 1896
$ws.Range("M132").Value = -5.521700000000237
$ws.Range("N132").Value = -6956

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27103.25
$ws.Range("I82").Value = 14206.5
$ws.Range("K82").Value = 14206.5
$ws.Range("M82").Value = -13823.5
$ws.Range("H85").Value = 27103.25
$ws.Range("I85").Value = 14206.5
$ws.Range("K85").Value = 14206.5
$ws.Range("M85").Value = -12880.5
$ws.Range("H99").Value = 5474.75
$ws.Range("J99").Value = 4332.6665
$ws.Range("L99").Value = 4332.6665
$ws.Range("N99").Value = -7328.6665
$ws.Range("H105").Value = 3242.8
$ws.Range("I105").Value = 3134.5
$ws.Range("K105").Value = 3134.5
$ws.Range("M105").Value = -1387.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2991.8235
$ws.Range("I31").Value = 2131.8
$ws.Range("K31").Value = 2131.8
$ws.Range("M31").Value = -1836.8
$ws.Range("H34").Value = 2991.8235
$ws.Range("I34").Value = 2131.8
$ws.Range("K34").Value = 2131.8
$ws.Range("M34").Value = -1929.8
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H92").Value = 43666.668
$ws.Range("J92").Value = 47500
$ws.Range("L92").Value = 47500
$ws.Range("N92").Value = -52492
$ws.Range("H105").Value = 1624.75
$ws.Range("I105").Value = 1250
$ws.Range("K105").Value = 1250
$ws.Range("M105").Value = 497
$ws.Range("H134").Value = 2536.8823
$ws.Range("I134").Value = 2093.5386
$ws.Range("J134").Value = 3977.75
$ws.Range("K134").Value = 6280.6158
$ws.Range("L134").Value = 11933.25
$ws.Range("M134").Value = -3745.6158
$ws.Range("N134").Value = -17003.25
$ws.Range("M45").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2535.25
$ws.Range("J55").Value = 3459.4
$ws.Range("L55").Value = 10378.2
$ws.Range("N55").Value = -10732.2
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H113").Value = 816.381
$ws.Range("I113").Value = 696.55554
$ws.Range("J113").Value = 906.25
$ws.Range("K113").Value = 2089.66662
$ws.Range("L113").Value = 2718.75
$ws.Range("M113").Value = 80.33338000000003
$ws.Range("N113").Value = -7058.75
$ws.Range("M58").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3160.2
$ws.Range("J41").Value = 3000
$ws.Range("L41").Value = 3000
$ws.Range("N41").Value = -3710
$ws.Range("H107").Value = 1550
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -3940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 68664
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H93").Value = 897
$ws.Range("I93").Value = 800
$ws.Range("K93").Value = 800
$ws.Range("M93").Value = 448
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H126").Value = 68664
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4200
$ws.Range("J96").Value = 4200
$ws.Range("L96").Value = 4200
$ws.Range("N96").Value = -6946
$ws.Range("H132").Value = 1389.5
$ws.Range("I132").Value = 1387.9
$ws.Range("J132").Value = 1397.5
$ws.Range("K132").Value = 4163.700000000001
$ws.Range("L132").Value = 4192.5
$ws.Range("M132").Value = -1633.700000000001
$ws.Range("N132").Value = -9252.5
